$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @('D2', '29.742.72'),
    @('E2', '  +1.90%  '),
    @('D3', '1.857.50'),
    @('E3', '  +1.46%  '),
    @('E4', '  -0.06%  '),
    @('D5', '244.67'),
    @('E5', '  +0.78%  '),
    @('D6', '0.6418'),
    @('E6', '  +3.51%  '),
    @('E7', '  -0.08%  '),
    @('D8', '47.41'),
    @('E8', '  +4.44%  '),
    @('E9', '  +2.33%  '),
    @('D10', '0.2980'),
    @('E10', '  +1.93%  '),
    @('E11', '  +5.37%  '),
    @('D12', '0.07668'),
    @('E12', '  +0.16%  '),
    @('D13', '1.872.96'),
    @('E13', '  +1.78%  '),
    @('E14', '  +1.37%  '),
    @('D15', '0.6914'),
    @('E15', '  +3.42%  '),
    @('D16', '83.85'),
    @('E16', '  +1.65%  '),
    @('D17', '0.000009840'),
    @('E17', '  +9.72%  '),
    @('D18', '6.095'),
    @('E18', '  +4.29%  '),
    @('D19', '29.751.56'),
    @('E19', '  +1.98%  '),
    @('D20', '2.115.72'),
    @('E20', '  +1.38%  '),
    @('D21', '235.98'),
    @('E21', '  +0.38%  '),
    @('E22', '  +1.65%  '),
    @('E23', '  -0.01%  '),
    @('D24', '7.515'),
    @('E24', '  +2.06%  '),
    @('D25', '1.000'),
    @('E25', '  +0.00%  '),
    @('D26', '158.91'),
    @('E26', '  +0.48%  '),
    @('D27', '0.1420'),
    @('E27', '  +1.82%  '),
    @('D28', '8.540'),
    @('D29', '17.91'),
    @('E29', '  +1.61%  '),
    @('D30', '0.06233'),
    @('E30', '  +7.13%  '),
    @('D31', '1.493'),
    @('E31', '  +0.43%  '),
    @('D32', '1.284'),
    @('E32', '  +6.20%  '),
    @('D33', '4.163'),
    @('E33', '  +1.90%  '),
    @('D34', '4.098'),
    @('D35', '1.898'),
    @('E35', '  +1.38%  '),
    @('E36', '  +2.84%  '),
    @('D37', '0.7273'),
    @('E37', '  +0.16%  '),
    @('E38', '  -0.24%  '),
    @('D39', '2.828'),
    @('E39', '  -1.07%  '),
    @('E40', '  +1.94%  '),
    @('D41', '1.203.21'),
    @('E41', '  -1.46%  '),
    @('D42', '0.9245'),
    @('E42', '  +1.94%  '),
    @('D43', '6.245'),
    @('E43', '  +0.02%  '),
    @('E44', '  +2.08%  '),
    @('D45', '0.9997'),
    @('E45', '  -0.03%  '),
    @('D46', '101.95'),
    @('E46', '  +0.25%  '),
    @('D47', '66.52'),
    @('E47', '  +1.62%  '),
    @('E48', '  +1.07%  '),
    @('B49', 'EnergySwap'),
    @('C49', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'),
    @('D49', '9.228'),
    @('B50', 'TheSandbox'),
    @('C50', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'),
    @('D50', '0.4060'),
    @('E50', '  +0.85%  '),
    @('E51', '  +0.84%  ')
)

foreach ($change in $changes) {
    $cellRef = $change[0]
    $newValue = $change[1]
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
}
